$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before F ("associate_charge"), shifting
#    service_charge..interview_end_date one column to the right (F -> G ... N -> O).
$ws.Columns("F:F").Insert()

# 2. Populate the new "associate_charge" column (header + value).
$ws.Range("F1").Value = "associate_charge"
$ws.Range("F2").Value = 1500

# 3. Give the new header cell the same look as the other "code-ish" headers
#    (e.g. interview_location, now shifted to I1) - Consolas / FFCE9178.
$null = $ws.Range("I1").Copy($ws.Range("F1"))

# Re-apply the text (Copy() also copies the shared-string value) - make sure
# F1 still reads "associate_charge" after the format copy.
$ws.Range("F1").Value = "associate_charge"

# 4. Drop the trailing "interview_end_date" column (now pushed out to O by
#    the insert above).
$ws.Columns("O:O").Delete()

# 5. The column insert does not carry the hyperlink anchor along with it, so
#    point it at the vendor_email value's new home (L2) by hand.
$ws.Range("K2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:vendor@yopmai.com")
$ws.Range("L2").Style = "Hyperlink"

# 6. Match the new column F width (~19.57 chars wide in the saved file; the
#    runtime quantizes ColumnWidth to 1/6ths, so 19.5 is the closest we can get).
$ws.Columns("F:F").ColumnWidth = 18.666666666666668

# 7. Reflect the new selection/scroll position left behind after the edit.
$ws.Range("O1:O1048576").Select()
